$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the Monday hours for the week commencing row 9 (A9 = 43171)
$ws.Range("B9").Value = 7

# Update the active selection to reflect where the user clicked next
$ws.Range("F11").Select()
